$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated simulation results (more games simulated; matrix re-normalized per row).
$ws.Range("B2").Value = 0.1877394636015326
$ws.Range("C2").Value = 0.5708812260536399
$ws.Range("J2").Value = 0.01915708812260536
$ws.Range("P2").Value = 0.1494252873563219
$ws.Range("S2").Value = 0.07279693486590039
$ws.Range("B3").Value = 0.006578947368421052
$ws.Range("C3").Value = 0.03289473684210526
$ws.Range("J3").Value = 0.006578947368421052
$ws.Range("P3").Value = 0.8026315789473685
$ws.Range("S3").Value = 0.1513157894736842
$ws.Range("J4").Value = 0.0975609756097561
$ws.Range("P4").Value = 0.5609756097560976
$ws.Range("S4").Value = 0.3414634146341464
$ws.Range("P5").Value = 0.75
$ws.Range("S5").Value = 0.25
$ws.Range("B6").Value = 0.09696969696969697
$ws.Range("D6").Value = 0.0303030303030303
$ws.Range("F6").Value = 0.03636363636363636
$ws.Range("J6").Value = 0.2545454545454545
$ws.Range("O6").Value = 0.02424242424242424
$ws.Range("Q6").Value = 0.07272727272727272
$ws.Range("R6").Value = 0.08484848484848485
$ws.Range("S6").Value = 0.4
$ws.Range("B7").Value = 0.1428571428571428
$ws.Range("D7").Value = 0.02285714285714286
$ws.Range("F7").Value = 0.05714285714285714
$ws.Range("J7").Value = 0.1714285714285714
$ws.Range("O7").Value = 0.01142857142857143
$ws.Range("Q7").Value = 0.1142857142857143
$ws.Range("R7").Value = 0.08
$ws.Range("S7").Value = 0.4
$ws.Range("B8").Value = 0.1160949868073879
$ws.Range("D8").Value = 0.01055408970976253
$ws.Range("E8").Value = 0.002638522427440633
$ws.Range("F8").Value = 0.06860158311345646
$ws.Range("J8").Value = 0.1345646437994723
$ws.Range("O8").Value = 0.01319261213720317
$ws.Range("Q8").Value = 0.1688654353562005
$ws.Range("R8").Value = 0.1029023746701847
$ws.Range("S8").Value = 0.3825857519788918
$ws.Range("B9").Value = 0.106508875739645
$ws.Range("D9").Value = 0.005917159763313609
$ws.Range("F9").Value = 0.04733727810650887
$ws.Range("J9").Value = 0.1183431952662722
$ws.Range("O9").Value = 0.01183431952662722
$ws.Range("Q9").Value = 0.1597633136094675
$ws.Range("R9").Value = 0.09467455621301775
$ws.Range("S9").Value = 0.4556213017751479
$ws.Range("B10").Value = 0.09532374100719425
$ws.Range("D10").Value = 0.02428057553956835
$ws.Range("E10").Value = 0.002697841726618705
$ws.Range("F10").Value = 0.06115107913669065
$ws.Range("J10").Value = 0.1348920863309352
$ws.Range("O10").Value = 0.01348920863309352
$ws.Range("Q10").Value = 0.1888489208633093
$ws.Range("R10").Value = 0.08992805755395683
$ws.Range("S10").Value = 0.3893884892086331
$ws.Range("G11").Value = 0.1098901098901099
$ws.Range("J11").Value = 0.08058608058608059
$ws.Range("K11").Value = 0.1758241758241758
$ws.Range("L11").Value = 0.6227106227106227
$ws.Range("S11").Value = 0.01098901098901099
$ws.Range("G12").Value = 0.7630057803468208
$ws.Range("J12").Value = 0.1734104046242775
$ws.Range("L12").Value = 0.0115606936416185
$ws.Range("S12").Value = 0.05202312138728324
$ws.Range("G13").Value = 0.7575757575757576
$ws.Range("J13").Value = 0.2424242424242424
$ws.Range("F15").Value = 0.01704545454545454
$ws.Range("H15").Value = 0.1704545454545454
$ws.Range("I15").Value = 0.08522727272727272
$ws.Range("J15").Value = 0.4147727272727273
$ws.Range("K15").Value = 0.04545454545454546
$ws.Range("M15").Value = 0.01136363636363636
$ws.Range("O15").Value = 0.02272727272727273
$ws.Range("S15").Value = 0.2329545454545454
$ws.Range("H16").Value = 0.1657458563535912
$ws.Range("I16").Value = 0.07734806629834254
$ws.Range("J16").Value = 0.3812154696132597
$ws.Range("K16").Value = 0.1933701657458564
$ws.Range("M16").Value = 0.03867403314917127
$ws.Range("O16").Value = 0.02209944751381215
$ws.Range("S16").Value = 0.1215469613259668
$ws.Range("F17").Value = 0.01201201201201201
$ws.Range("H17").Value = 0.1711711711711712
$ws.Range("I17").Value = 0.09309309309309309
$ws.Range("J17").Value = 0.3993993993993994
$ws.Range("K17").Value = 0.1141141141141141
$ws.Range("M17").Value = 0.01501501501501501
$ws.Range("O17").Value = 0.07807807807807808
$ws.Range("S17").Value = 0.1171171171171171
$ws.Range("F18").Value = 0.03314917127071823
$ws.Range("H18").Value = 0.143646408839779
$ws.Range("I18").Value = 0.09944751381215469
$ws.Range("J18").Value = 0.4198895027624309
$ws.Range("K18").Value = 0.1160220994475138
$ws.Range("O18").Value = 0.0718232044198895
$ws.Range("S18").Value = 0.1160220994475138
$ws.Range("F19").Value = 0.01078167115902965
$ws.Range("H19").Value = 0.2129380053908356
$ws.Range("I19").Value = 0.08265947888589398
$ws.Range("J19").Value = 0.3719676549865229
$ws.Range("K19").Value = 0.1087151841868823
$ws.Range("M19").Value = 0.01886792452830189
$ws.Range("N19").Value = 0.0008984725965858042
$ws.Range("O19").Value = 0.06828391734052111
$ws.Range("S19").Value = 0.1248876909254268

Write-Output "Updated 107 cells on $($ws.Name)"
